$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "about_me"
$ws.Range("D2").Value = "About Me"

$ws.Range("D3").Select()
